# Auto-generated edit script: updates crypto price/volume table cells
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.142.47"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.654.52"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.80"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5237"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06358"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.50"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07701"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("D13").Value = "1.736.39"
$ws.Range("E13").Value = "  +5.04%  "
$ws.Range("D14").Value = "1.882.16"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5617"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "0.0₅8202"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.54"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "26.129.96"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.653"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.33"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.959"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.23"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1196"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.265"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.96"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.515"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05477"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.271"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.465"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.380"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.559"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9540"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5664"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01581"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.868"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8360"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "1.029.07"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.21"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "1.792.85"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.87"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "0.0₈109"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.016"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4339"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05196"
$ws.Range("E51").Value = "  -4.05%  "
